$d = $word.ActiveDocument

# Helper characters
$rsq = [char]0x2019   # right single quotation mark (')
$ndash = [char]0x2013 # en dash (-)

function Find-Replace($searchText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: text not found: $searchText"
    }
    return $ok
}

# ---------------------------------------------------------------------------
# 1. "Led development of the Department's policy bot and first retrieval-
#    augmented generation implementation..." bullet -> rewritten single run
# ---------------------------------------------------------------------------
$old1 = "Led development of the Department" + $rsq + "s policy bot and first retrieval-augmented generation implementation in the StateChat suite, spanning 25,000+ pages of policy. Implemented a variable-length semantic chunking process and custom-built hybrid-search process based on a novel evaluation framework (incorporating metrics for answer accuracy, citation accuracy, and inference speed)"
$new1 = "Led development of the Department" + $rsq + "s first retrieval-augmented generation (RAG) system, integrating 25,000+ pages of policy using a variable-length semantic chunking and custom hybrid search approach. Designed a novel evaluation framework incorporating answer and citation accuracy and inference speed"
Find-Replace $old1 $new1

# ---------------------------------------------------------------------------
# 2. "Developed a tool to enable foreign service officers..." bullet ->
#    rewritten, split across 3 runs (same formatting) per target OOXML
# ---------------------------------------------------------------------------
$old2 = "Developed a tool to enable foreign service officers to ask questions before bidding on their next assignment" + [char]32 + $ndash + "the Department" + $rsq + "s first chain-of-thought function-calling implementation, enabling the tool to retrieve information from multiple sources and conduct variable number searches based on user queries"
$part2a = "Built a chain-of-thought function-calling system "
$part2b = "(the Department" + $rsq + "s first) "
$part2c = "to assist foreign service officers in assignment bidding, enabling multi-source retrieval and adaptive query searches based on user intent"
$new2 = $part2a + $part2b + $part2c
Find-Replace $old2 $new2

# Split the merged run into 3 runs (same rPr) by toggling bold off/on which
# forces Word to re-split runs at the toggle boundaries without altering the
# visible formatting (bCs was already absent of plain bold).
$rngSplit = $d.Content
$rngSplit.Find.Execute($new2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rngSplit.Find.Found) {
    $len2a = $part2a.Length
    $len2b = $part2b.Length
    $sub2a = $d.Range($rngSplit.Start, $rngSplit.Start + $len2a)
    $sub2a.Font.Bold = 1
    $sub2a.Font.Bold = 0
    $sub2b = $d.Range($rngSplit.Start + $len2a, $rngSplit.Start + $len2a + $len2b)
    $sub2b.Font.Bold = 1
    $sub2b.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# 3. "Designed trainings, promotional materials..." bullet -> split into two
#    separate bullets: "Managed versioning..." and "Drove adoption..."
# ---------------------------------------------------------------------------
$old3 = "Designed trainings, promotional materials, and user testing artifacts for enterprise products, and led dozens of trainings of up to 450 people for different user groups across the globe"
$part3a = "Managed versioning, testing, and deployment to production, ensuring reliability and scalability"
$part3b = "Drove adoption through data-informed training strategies, creating instructional materials and leading global training sessions for up to 450 users"
Find-Replace $old3 $part3a

$rng3 = $d.Content
$rng3.Find.Execute($part3a, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng3.Find.Found) {
    $rng3.Collapse(0)   # wdCollapseEnd
    $rng3.InsertParagraphAfter()
    $insertRng = $d.Range($rng3.End + 1, $rng3.End + 1)
    $insertRng.InsertAfter($part3b)
}

# ---------------------------------------------------------------------------
# 4. "Researched and compared multiple open-source LLMs..." -> "Compared..."
# ---------------------------------------------------------------------------
$old4 = "Researched and compared multiple open-source LLMs for accuracy and completeness for extraction tasks; selected a model that had not been released at the start of the campaign. Optimized hyperparameters for performance"
$new4 = "Compared multiple open-source LLMs for accuracy and completeness for extraction tasks; selected a model that had not been released at the start of the campaign. Optimized hyperparameters for performance"
Find-Replace $old4 $new4

# ---------------------------------------------------------------------------
# 5. "Designed data ingestion process using Microsoft Azure Cognitive
#    Services, including Document Intelligence, Vision, and Translator" ->
#    "...using Azure Cognitive Services (Document Intelligence, Vision,
#    Translator)"
# ---------------------------------------------------------------------------
$old5 = "Designed data ingestion process using Microsoft Azure Cognitive Services, including Document Intelligence, Vision, and Translator"
$new5 = "Designed data ingestion process using Azure Cognitive Services (Document Intelligence, Vision, Translator)"
Find-Replace $old5 $new5

# ---------------------------------------------------------------------------
# 6. Move w:lastRenderedPageBreak from "Strategic Competition Funding: " run
#    to the "Family Advocacy Program..." run.
# ---------------------------------------------------------------------------
# 6a. Remove from "Strategic Competition Funding: " by self-replacing the
#     text, which forces Word to regenerate that run without the stale
#     rendering-pagination hint.
Find-Replace "Strategic Competition Funding: " "Strategic Competition Funding: "

# 6b. Re-insert it immediately before "Family Advocacy Program..." by
#     replacing that run's XML with an equivalent run that also carries the
#     lastRenderedPageBreak marker.
$rng6 = $d.Content
$rng6.Find.Execute("Family Advocacy Program, United States Department of Defense", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng6.Find.Found) {
    $fragment6 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:eastAsia="Garamond" w:hAnsi="Garamond" w:cs="Garamond"/><w:bCs/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>Family Advocacy Program, United States Department of Defense</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng6.InsertXML($fragment6)
}

# ---------------------------------------------------------------------------
# 7. Spelling fixes: Pytorch -> PyTorch, Pyspark -> PySpark
# ---------------------------------------------------------------------------
Find-Replace "Pytorch" "PyTorch"
Find-Replace "Pyspark" "PySpark"
